$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.3464964993005633, 86.29678392075563, 3.082599426703578, 246.9852506941017, 336.7111305408614)
    3 = @(1.505614041169197, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 6.741336633845642)
    4 = @(1.505614041169197, 0.05231270169004087, 0.7127328510149897, 0.4998867070740569, 2.770546300948285)
    5 = @(3.182878228561681, 1.65323645889881, 16.98373111632243, 0.4998867070740569, 22.31973251085698)
    6 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    7 = @(0.3464964993005633, 9.226618575922256, 3.082599426703578, 6.48142807727062, 19.13714257919702)
    8 = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 5.488907176552729)
    9 = @(0.7287194209349384, 0.004309184025731883, 0.1529057820181812, 0.4998867070740569, 1.385821094052908)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 7).Value = $vals[4]
}
